$p = $ppt.ActivePresentation

# 1) "Atualização no Slide": remove the placeholder slide that only
#    contains "?" (currently slide #16 in the deck). The following
#    slides (Institucional & Dashboard, JFrame, Cliente em Maquina
#    Virtual, Obrigado!) all shift up one position and the deck ends
#    up with 19 slides instead of 20.
$p.Slides.Item(16).Delete()

# 2) "correção gramatical": tidy up the stray trailing empty run
#    (endParaRPr) left after the "Web Designer" text on slide 2's
#    second subtitle placeholder. Re-typing the same text removes the
#    redundant empty end-of-paragraph run properties while keeping the
#    run's own formatting (size 28, green fill) untouched.
$s2 = $p.Slides.Item(2)
$webDesignerShape = $s2.Shapes.Item(6)
$tr = $webDesignerShape.TextFrame.TextRange
$tr.Delete()
$tr.Text = "Web Designer"
